# Macroferia Regional de Talca - Coliflor: insert a new weekly price
# observation row. This shifts every existing record from row 308
# onward down by one row (308->309, ..., 372->373) and fills the
# freed-up row 308 with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 308, pushing rows 308..372 down to 309..373.
$ws.Rows("308").Insert()

# Populate the new row 308 with the new weekly record.
$ws.Cells.Item(308, 1).Value = 5
$ws.Cells.Item(308, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(308, 3).Value = "Maule"
$ws.Cells.Item(308, 4).Value = 44889
$ws.Cells.Item(308, 5).Value = 7
$ws.Cells.Item(308, 6).Value = 100112008
$ws.Cells.Item(308, 7).Value = "Coliflor"
$ws.Cells.Item(308, 8).Value = "Sin especificar"
$ws.Cells.Item(308, 9).Value = "Primera"
$ws.Cells.Item(308, 10).Value = 5000
$ws.Cells.Item(308, 11).Value = 600
$ws.Cells.Item(308, 12).Value = 600
$ws.Cells.Item(308, 13).Value = 600
$ws.Cells.Item(308, 14).Value = "`$/unidad"
$ws.Cells.Item(308, 15).Value = "Región del Maule"
$ws.Cells.Item(308, 16).Value = 600
$ws.Cells.Item(308, 17).Value = 1
$ws.Cells.Item(308, 18).Value = "Hortaliza"
